# Update cryptos list - GitHub Actions data refresh (Sun May 28 16:25:11 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that Excel would otherwise auto-convert to a number
# (e.g. "1.012", "0.4768") as plain text, matching the source file's
# inline-string cells, without leaving a lingering quote-prefixed style on
# the cell.
function Set-TextValue {
    param($addr, $value)
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.504.57"
$ws.Range("E2").Value = "  +1.94%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.858.95"

# Row 4 - TetherUSD
Set-TextValue "D4" "1.012"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5 - BNB
Set-TextValue "D5" "311.68"
$ws.Range("E5").Value = "  +0.79%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.31%  "

# Row 7 - XRP
Set-TextValue "D7" "0.4768"
$ws.Range("E7").Value = "  +0.08%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3793"
$ws.Range("E8").Value = "  +3.22%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.07316"
$ws.Range("E9").Value = "  +1.32%  "

# Row 10 - Polygon
Set-TextValue "D10" "0.9289"
$ws.Range("E10").Value = "  -0.34%  "

# Row 11 - Solana
Set-TextValue "D11" "20.68"
$ws.Range("E11").Value = "  +4.15%  "

# Row 12 - TRON
Set-TextValue "D12" "0.07777"
$ws.Range("E12").Value = "  +0.36%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.893.15"
$ws.Range("E13").Value = "  +2.53%  "

# Row 14 - Polkadot
Set-TextValue "D14" "5.432"
$ws.Range("E14").Value = "  +0.88%  "

# Row 15 - Chainlink
Set-TextValue "D15" "6.556"
$ws.Range("E15").Value = "  +1.31%  "

# Row 16 - Litecoin
Set-TextValue "D16" "90.08"
$ws.Range("E16").Value = "  +1.31%  "

# Row 17 - BinanceUSD
Set-TextValue "D17" "1.013"
$ws.Range("E17").Value = "  -0.39%  "

# Row 18 - ShibaInu
Set-TextValue "D18" "0.000008792"
$ws.Range("E18").Value = "  +1.58%  "

# Row 19 - Dai
Set-TextValue "D19" "1.010"
$ws.Range("E19").Value = "  -0.28%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "27.550.90"
$ws.Range("E20").Value = "  +1.94%  "

# Row 21 - Avalanche
Set-TextValue "D21" "14.64"
$ws.Range("E21").Value = "  +0.64%  "

# Row 22 - Uniswap
Set-TextValue "D22" "5.089"
$ws.Range("E22").Value = "  +0.56%  "

# Row 23 - Cosmos
Set-TextValue "D23" "10.69"
$ws.Range("E23").Value = "  +0.48%  "

# Row 24 - Toncoin
Set-TextValue "D24" "1.935"
$ws.Range("E24").Value = "  +0.26%  "

# Row 25 - Monero
Set-TextValue "D25" "155.81"
$ws.Range("E25").Value = "  +1.90%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "18.46"
$ws.Range("E26").Value = "  +1.19%  "

# Row 27 - LidoDAOToken
Set-TextValue "D27" "2.004"
$ws.Range("E27").Value = "  +0.37%  "

# Row 28 - BitcoinCash
Set-TextValue "D28" "115.29"
$ws.Range("E28").Value = "  +0.72%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue "D29" "4.940"
$ws.Range("E29").Value = "  -0.59%  "

# Row 30 - Stellar
Set-TextValue "D30" "0.08860"
$ws.Range("E30").Value = "  -0.08%  "

# Row 31 - HuobiToken
Set-TextValue "D31" "3.324"
$ws.Range("E31").Value = "  +0.39%  "

# Row 32 - ARBITRUM
Set-TextValue "D32" "1.201"
$ws.Range("E32").Value = "  +1.76%  "

# Row 33 - ImmutableX
Set-TextValue "D33" "0.7505"
$ws.Range("E33").Value = "  +1.78%  "

# Row 34 - Filecoin
Set-TextValue "D34" "4.576"
$ws.Range("E34").Value = "  +1.59%  "

# Row 35 - RenderToken
Set-TextValue "D35" "2.687"
$ws.Range("E35").Value = "  +0.03%  "

# Rows 36/37 swap places: VeChain <-> TrustWalletToken, with updated values
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D36" "1.121"
$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D37" "0.02034"
$ws.Range("E37").Value = "  +3.27%  "

# Row 38 - TheSandbox
$ws.Range("E38").Value = "  +6.49%  "

# Row 39 - Hedera
Set-TextValue "D39" "0.05316"
$ws.Range("E39").Value = "  +0.87%  "

# Row 40 - MXToken
$ws.Range("E40").Value = "  +0.10%  "

# Row 41 - FraxShare
Set-TextValue "D41" "6.996"
$ws.Range("E41").Value = "  -0.41%  "

# Rows 42/43 swap places: Aptos <-> Algorand, with updated values
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D42" "0.1520"
$ws.Range("E42").Value = "  +0.50%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D43" "8.450"
$ws.Range("E43").Value = "  +2.02%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "10.70"
$ws.Range("E44").Value = "  +1.30%  "

# Row 45 - Decentraland
Set-TextValue "D45" "0.4849"
$ws.Range("E45").Value = "  +2.37%  "

# Row 46 - PaxDollar
$ws.Range("E46").Value = "  -0.44%  "

# Row 47 - Quant
Set-TextValue "D47" "104.15"
$ws.Range("E47").Value = "  +2.34%  "

# Row 48 - NEARProtocol
Set-TextValue "D48" "1.661"
$ws.Range("E48").Value = "  +3.11%  "

# Row 49 - Aave
Set-TextValue "D49" "67.39"
$ws.Range("E49").Value = "  +2.57%  "

# Row 50 - Cronos
Set-TextValue "D50" "0.06099"
$ws.Range("E50").Value = "  +0.62%  "

# Row 51 - EOS
Set-TextValue "D51" "0.9085"
$ws.Range("E51").Value = "  +2.01%  "
